$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    if ($val -eq "" -or $val -match "^\d{4}-\d{2}-\d{2}$") {
        $ws.Range($addr).Formula = "'" + $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}

# --- Row 3 ---
$ws.Range('A3').Value = 130853761
$ws.Range('B3').Value = 79244
Set-TextValue $ws 'D3' 'NT'
$ws.Range('E3').Value = 230405
Set-TextValue $ws 'F3' 'Garnlav (ssp. sarmentosa)'
Set-TextValue $ws 'G3' 'Alectoria sarmentosa subsp. sarmentosa'
Set-TextValue $ws 'H3' '(Ach.) Ach.'
Set-TextValue $ws 'I3' ''
Set-TextValue $ws 'P3' 'Djupbäcken, Djupbäcken, Jmt'
$ws.Range('Q3').Value = 442771
$ws.Range('R3').Value = 7039709
$ws.Range('S3').Value = 20
Set-TextValue $ws 'T3' 'Jämtland'
Set-TextValue $ws 'U3' 'Krokom'
Set-TextValue $ws 'V3' 'Jämtland'
Set-TextValue $ws 'W3' 'Offerdal'
Set-TextValue $ws 'Y3' '2026-01-24'
Set-TextValue $ws 'Z3' '11:05'
Set-TextValue $ws 'AA3' '2026-01-24'
Set-TextValue $ws 'AB3' '11:05'
$ws.Range('AD3').Value = $false
$ws.Range('AE3').Value = $false
$ws.Range('AG3').Value = $false
Set-TextValue $ws 'AT3' ''
Set-TextValue $ws 'AW3' 'Maria Danvind'
Set-TextValue $ws 'AX3' 'Maria Danvind'
Set-TextValue $ws 'AY3' ''
$ws.Range('AC3').Value = ''

# --- Row 4 ---
$ws.Range('A4').Value = 130861152
$ws.Range('B4').Value = 91804
Set-TextValue $ws 'D4' 'NT'
$ws.Range('E4').Value = 1108
Set-TextValue $ws 'F4' 'Harticka'
Set-TextValue $ws 'G4' 'Pelloporus leporinus'
Set-TextValue $ws 'H4' '(Fr.) Krieglst.'
Set-TextValue $ws 'I4' ''
Set-TextValue $ws 'P4' 'Djupbäcken, Jmt'
$ws.Range('Q4').Value = 442868
$ws.Range('R4').Value = 7039767
$ws.Range('S4').Value = 10
Set-TextValue $ws 'T4' 'Jämtland'
Set-TextValue $ws 'U4' 'Krokom'
Set-TextValue $ws 'V4' 'Jämtland'
Set-TextValue $ws 'W4' 'Offerdal'
Set-TextValue $ws 'Y4' '2026-01-24'
Set-TextValue $ws 'AA4' '2026-01-24'
Set-TextValue $ws 'AC4' 'I stående levande gran med full längd.'
$ws.Range('AD4').Value = $false
$ws.Range('AE4').Value = $false
$ws.Range('AG4').Value = $false
Set-TextValue $ws 'AT4' ''
Set-TextValue $ws 'AW4' 'Kristian Zackrisson'
Set-TextValue $ws 'AX4' 'Kristian Zackrisson'
Set-TextValue $ws 'AY4' ''
$ws.Range('AB4').Value = ''
$ws.Range('Z4').Value = ''

# --- Row 10 ---
$ws.Range('A10').Value = 130861151
$ws.Range('B10').Value = 57884
Set-TextValue $ws 'D10' 'NT'
$ws.Range('E10').Value = 100109
Set-TextValue $ws 'F10' 'Tretåig hackspett'
Set-TextValue $ws 'G10' 'Picoides tridactylus'
Set-TextValue $ws 'H10' '(Linnaeus, 1758)'
Set-TextValue $ws 'I10' ''
Set-TextValue $ws 'K10' ''
Set-TextValue $ws 'L10' ''
Set-TextValue $ws 'M10' 'äldre spår'
Set-TextValue $ws 'N10' ''
Set-TextValue $ws 'P10' 'Djupbäcken, Jmt'
$ws.Range('Q10').Value = 442749
$ws.Range('R10').Value = 7039568
$ws.Range('S10').Value = 10
Set-TextValue $ws 'T10' 'Jämtland'
Set-TextValue $ws 'U10' 'Krokom'
Set-TextValue $ws 'V10' 'Jämtland'
Set-TextValue $ws 'W10' 'Offerdal'
Set-TextValue $ws 'Y10' '2026-01-24'
Set-TextValue $ws 'AA10' '2026-01-24'
Set-TextValue $ws 'AC10' 'Ringhack, äldre, enstaka på en gran. Mycket högt livsmiljövärde för tretåig hackspett kring fyndplatsen.'
$ws.Range('AD10').Value = $false
$ws.Range('AE10').Value = $false
$ws.Range('AG10').Value = $false
Set-TextValue $ws 'AH10' 'Granskog'
Set-TextValue $ws 'AJ10' 'gran'
Set-TextValue $ws 'AK10' 'Picea abies'
Set-TextValue $ws 'AM10' 'Trädstam på levande träd'
Set-TextValue $ws 'AO10' 'Stem on living tree # Picea abies'
Set-TextValue $ws 'AT10' ''
Set-TextValue $ws 'AW10' 'Kristian Zackrisson'
Set-TextValue $ws 'AX10' 'Kristian Zackrisson'
Set-TextValue $ws 'AY10' ''

# --- Row 11 ---
$ws.Range('A11').Value = 130861155
$ws.Range('B11').Value = 79243
Set-TextValue $ws 'D11' 'NT'
$ws.Range('E11').Value = 6425
Set-TextValue $ws 'F11' 'Garnlav'
Set-TextValue $ws 'G11' 'Alectoria sarmentosa'
Set-TextValue $ws 'H11' '(Ach.) Ach.'
Set-TextValue $ws 'I11' ''
Set-TextValue $ws 'P11' 'Djupbäcken, Jmt'
$ws.Range('Q11').Value = 442870
$ws.Range('R11').Value = 7039632
$ws.Range('S11').Value = 10
Set-TextValue $ws 'T11' 'Jämtland'
Set-TextValue $ws 'U11' 'Krokom'
Set-TextValue $ws 'V11' 'Jämtland'
Set-TextValue $ws 'W11' 'Offerdal'
Set-TextValue $ws 'Y11' '2026-01-24'
Set-TextValue $ws 'AA11' '2026-01-24'
Set-TextValue $ws 'AC11' 'På död undertryck gran.'
$ws.Range('AD11').Value = $false
$ws.Range('AE11').Value = $false
$ws.Range('AG11').Value = $false
Set-TextValue $ws 'AT11' ''
Set-TextValue $ws 'AW11' 'Kristian Zackrisson'
Set-TextValue $ws 'AX11' 'Kristian Zackrisson'
Set-TextValue $ws 'AY11' ''
$ws.Range('AH11').Value = ''
$ws.Range('AJ11').Value = ''
$ws.Range('AK11').Value = ''
$ws.Range('AM11').Value = ''
$ws.Range('AO11').Value = ''
$ws.Range('K11').Value = ''
$ws.Range('L11').Value = ''
$ws.Range('M11').Value = ''
$ws.Range('N11').Value = ''
